# Update "想去人数" (want-to-go count) values in column F for the rows whose
# events had updated attendance counts. The two sheets "展览" and "全部类型"
# contain the same underlying data set, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 291
    6  = 314
    7  = 9247
    9  = 79
    12 = 10
    17 = 270
    18 = 771
    19 = 43
    20 = 91
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
